$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 12.52413133333333
$ws.Range("H2").Value = 37.572394
$ws.Range("I2").Value = 0.09718402715578596
$ws.Range("J2").Value = 0.1008592412859651
$ws.Range("M2").Value = 29.785352
$ws.Range("N2").Value = 89.356056
$ws.Range("O2").Value = 0.7923195065866085
$ws.Range("P2").Value = 0.7947519366640845
$ws.Range("Q2").Value = 373.0356602575626
$ws.Range("R2").Value = 3357.320942318063
$ws.Range("S2").Value = 0.0770008004441719
$ws.Range("T2").Value = 0.08015807734249097

$ws.Range("G3").Value = 12.52413133333333
$ws.Range("H3").Value = 37.572394
$ws.Range("I3").Value = 0.09718402715578596
$ws.Range("J3").Value = 0.1008592412859651
$ws.Range("O3").Value = 0.1171985110386058
$ws.Range("P3").Value = 0.1175583118271966
$ws.Range("Q3").Value = 55.17878025600599
$ws.Range("R3").Value = 496.6090223040539
$ws.Range("S3").Value = 0.01138982327939355
$ws.Range("T3").Value = 0.01185684213774995

$ws.Range("G4").Value = 12.52413133333333
$ws.Range("H4").Value = 37.572394
$ws.Range("I4").Value = 0.09718402715578596
$ws.Range("J4").Value = 0.1008592412859651
$ws.Range("M4").Value = 1.362560333333333
$ws.Range("N4").Value = 4.087681
$ws.Range("O4").Value = 0.03624543805965938
$ws.Range("P4").Value = 0.03635671197501131
$ws.Range("Q4").Value = 17.06488456425711
$ws.Range("R4").Value = 153.583961078314
$ws.Range("S4").Value = 0.003522477636663295
$ws.Range("T4").Value = 0.003666910385452004

$ws.Range("G5").Value = 12.52413133333333
$ws.Range("H5").Value = 37.572394
$ws.Range("I5").Value = 0.09718402715578596
$ws.Range("J5").Value = 0.1008592412859651
$ws.Range("M5").Value = 0.3451695
$ws.Range("N5").Value = 0.690339
$ws.Range("O5").Value = 0.00918184643004207
$ws.Range("P5").Value = 0.006140023203404898
$ws.Range("Q5").Value = 4.322948150260999
$ws.Range("R5").Value = 25.937688901566
$ws.Range("S5").Value = 0.0008923288127974649
$ws.Range("T5").Value = 0.0006192780817736392

$ws.Range("G6").Value = 12.52413133333333
$ws.Range("H6").Value = 37.572394
$ws.Range("I6").Value = 0.09718402715578596
$ws.Range("J6").Value = 0.1008592412859651
$ws.Range("M6").Value = 1.693723333333333
$ws.Range("N6").Value = 5.08117
$ws.Range("O6").Value = 0.04505469788508434
$ws.Range("P6").Value = 0.04519301633030275
$ws.Range("Q6").Value = 21.21241346899778
$ws.Range("R6").Value = 190.91172122098
$ws.Range("S6").Value = 0.004378596982759769
$ws.Range("T6").Value = 0.004558133338498567

$ws.Range("I7").Value = 0.1842225641940495
$ws.Range("J7").Value = 0.1911893198517306
$ws.Range("M7").Value = 29.785352
$ws.Range("N7").Value = 89.356056
$ws.Range("O7").Value = 0.7923195065866085
$ws.Range("P7").Value = 0.7947519366640845
$ws.Range("Q7").Value = 707.1284024720213
$ws.Range("R7").Value = 6364.155622248191
$ws.Range("S7").Value = 0.1459631311643491
$ws.Range("T7").Value = 0.151948082221652

$ws.Range("I8").Value = 0.1842225641940495
$ws.Range("J8").Value = 0.1911893198517306
$ws.Range("O8").Value = 0.1171985110386058
$ws.Range("P8").Value = 0.1175583118271966
$ws.Range("S8").Value = 0.02159061022325658
$ws.Range("T8").Value = 0.02247589368115938

$ws.Range("I9").Value = 0.1842225641940495
$ws.Range("J9").Value = 0.1911893198517306
$ws.Range("M9").Value = 1.362560333333333
$ws.Range("N9").Value = 4.087681
$ws.Range("O9").Value = 0.03624543805965938
$ws.Range("P9").Value = 0.03635671197501131
$ws.Range("Q9").Value = 32.34828689557689
$ws.Range("R9").Value = 291.134582060192
$ws.Range("S9").Value = 0.006677227539687044
$ws.Range("T9").Value = 0.006951015034547683

$ws.Range("I10").Value = 0.1842225641940495
$ws.Range("J10").Value = 0.1911893198517306
$ws.Range("M10").Value = 0.3451695
$ws.Range("N10").Value = 0.690339
$ws.Range("O10").Value = 0.00918184643004207
$ws.Range("P10").Value = 0.006140023203404898
$ws.Range("Q10").Value = 8.194603747407999
$ws.Range("R10").Value = 49.167622484448
$ws.Range("S10").Value = 0.001691503293378329
$ws.Range("T10").Value = 0.001173906860132827

$ws.Range("I11").Value = 0.1842225641940495
$ws.Range("J11").Value = 0.1911893198517306
$ws.Range("M11").Value = 1.693723333333333
$ws.Range("N11").Value = 5.08117
$ws.Range("O11").Value = 0.04505469788508434
$ws.Range("P11").Value = 0.04519301633030275
$ws.Range("Q11").Value = 40.21036497838222
$ws.Range("R11").Value = 361.89328480544
$ws.Range("S11").Value = 0.008300091973378455
$ws.Range("T11").Value = 0.008640422054238737

$ws.Range("G12").Value = 41.01852933333333
$ws.Range("H12").Value = 123.055588
$ws.Range("I12").Value = 0.3182932023406124
$ws.Range("J12").Value = 0.3303301152883236
$ws.Range("M12").Value = 29.785352
$ws.Range("N12").Value = 89.356056
$ws.Range("O12").Value = 0.7923195065866085
$ws.Range("P12").Value = 0.7947519366640845
$ws.Range("Q12").Value = 1221.751334715659
$ws.Range("R12").Value = 10995.76201244093
$ws.Range("S12").Value = 0.2521899130283856
$ws.Range("T12").Value = 0.2625304988638655

$ws.Range("G13").Value = 41.01852933333333
$ws.Range("H13").Value = 123.055588
$ws.Range("I13").Value = 0.3182932023406124
$ws.Range("J13").Value = 0.3303301152883236
$ws.Range("O13").Value = 0.1171985110386058
$ws.Range("P13").Value = 0.1175583118271966
$ws.Range("Q13").Value = 180.719313481212
$ws.Range("R13").Value = 1626.473821330908
$ws.Range("S13").Value = 0.03730348938802946
$ws.Range("T13").Value = 0.03883305069897854

$ws.Range("G14").Value = 41.01852933333333
$ws.Range("H14").Value = 123.055588
$ws.Range("I14").Value = 0.3182932023406124
$ws.Range("J14").Value = 0.3303301152883236
$ws.Range("M14").Value = 1.362560333333333
$ws.Range("N14").Value = 4.087681
$ws.Range("O14").Value = 0.03624543805965938
$ws.Range("P14").Value = 0.03635671197501131
$ws.Range("Q14").Value = 55.89022100126978
$ws.Range("R14").Value = 503.011989011428
$ws.Range("S14").Value = 0.01153667655024729
$ws.Range("T14").Value = 0.01200971685820986

$ws.Range("G15").Value = 41.01852933333333
$ws.Range("H15").Value = 123.055588
$ws.Range("I15").Value = 0.3182932023406124
$ws.Range("J15").Value = 0.3303301152883236
$ws.Range("M15").Value = 0.3451695
$ws.Range("N15").Value = 0.690339
$ws.Range("O15").Value = 0.00918184643004207
$ws.Range("P15").Value = 0.006140023203404898
$ws.Range("Q15").Value = 14.158345260722
$ws.Range("R15").Value = 84.95007156433201
$ws.Range("S15").Value = 0.00292251930361781
$ws.Range("T15").Value = 0.002028234572653722

$ws.Range("G16").Value = 41.01852933333333
$ws.Range("H16").Value = 123.055588
$ws.Range("I16").Value = 0.3182932023406124
$ws.Range("J16").Value = 0.3303301152883236
$ws.Range("M16").Value = 1.693723333333333
$ws.Range("N16").Value = 5.08117
$ws.Range("O16").Value = 0.04505469788508434
$ws.Range("P16").Value = 0.04519301633030275
$ws.Range("Q16").Value = 69.47404023088446
$ws.Range("R16").Value = 625.26636207796
$ws.Range("S16").Value = 0.01434060407033231
$ws.Range("T16").Value = 0.014928614294616

$ws.Range("G17").Value = 14.087727
$ws.Range("H17").Value = 28.175454
$ws.Range("I17").Value = 0.1093171260259301
$ws.Range("J17").Value = 0.07563411885139956
$ws.Range("M17").Value = 29.785352
$ws.Range("N17").Value = 89.356056
$ws.Range("O17").Value = 0.7923195065866085
$ws.Range("P17").Value = 0.7947519366640845
$ws.Range("Q17").Value = 419.607907574904
$ws.Range("R17").Value = 2517.647445449424
$ws.Range("S17").Value = 0.086614091354331
$ws.Range("T17").Value = 0.06011036243503134

$ws.Range("G18").Value = 14.087727
$ws.Range("H18").Value = 28.175454
$ws.Range("I18").Value = 0.1093171260259301
$ws.Range("J18").Value = 0.07563411885139956
$ws.Range("O18").Value = 0.1171985110386058
$ws.Range("P18").Value = 0.1175583118271966
$ws.Range("Q18").Value = 62.067665353419
$ws.Range("R18").Value = 372.405992120514
$ws.Range("S18").Value = 0.01281180440125863
$ws.Range("T18").Value = 0.008891419328708077

$ws.Range("G19").Value = 14.087727
$ws.Range("H19").Value = 28.175454
$ws.Range("I19").Value = 0.1093171260259301
$ws.Range("J19").Value = 0.07563411885139956
$ws.Range("M19").Value = 1.362560333333333
$ws.Range("N19").Value = 4.087681
$ws.Range("O19").Value = 0.03624543805965938
$ws.Range("P19").Value = 0.03635671197501131
$ws.Range("Q19").Value = 19.195377997029
$ws.Range("R19").Value = 115.172267982174
$ws.Range("S19").Value = 0.003962247120232826
$ws.Range("T19").Value = 0.002749807874564107

$ws.Range("G20").Value = 14.087727
$ws.Range("H20").Value = 28.175454
$ws.Range("I20").Value = 0.1093171260259301
$ws.Range("J20").Value = 0.07563411885139956
$ws.Range("M20").Value = 0.3451695
$ws.Range("N20").Value = 0.690339
$ws.Range("O20").Value = 0.00918184643004207
$ws.Range("P20").Value = 0.006140023203404898
$ws.Range("Q20").Value = 4.8626536847265
$ws.Range("R20").Value = 19.450614738906
$ws.Range("S20").Value = 0.001003733063343645
$ws.Range("T20").Value = 0.0004643952447166771

$ws.Range("G21").Value = 14.087727
$ws.Range("H21").Value = 28.175454
$ws.Range("I21").Value = 0.1093171260259301
$ws.Range("J21").Value = 0.07563411885139956
$ws.Range("M21").Value = 1.693723333333333
$ws.Range("N21").Value = 5.08117
$ws.Range("O21").Value = 0.04505469788508434
$ws.Range("P21").Value = 0.04519301633030275
$ws.Range("Q21").Value = 23.86071193353
$ws.Range("R21").Value = 143.16427160118
$ws.Range("S21").Value = 0.004925250086763969
$ws.Range("T21").Value = 0.003418133968379359

$ws.Range("G22").Value = 37.49906666666667
$ws.Range("H22").Value = 112.4972
$ws.Range("I22").Value = 0.2909830802836222
$ws.Range("J22").Value = 0.3019872047225811
$ws.Range("M22").Value = 29.785352
$ws.Range("N22").Value = 89.356056
$ws.Range("O22").Value = 0.7923195065866085
$ws.Range("P22").Value = 0.7947519366640845
$ws.Range("Q22").Value = 1116.922900338134
$ws.Range("R22").Value = 10052.3061030432
$ws.Range("S22").Value = 0.2305515705953711
$ws.Range("T22").Value = 0.2400049158010447

$ws.Range("G23").Value = 37.49906666666667
$ws.Range("H23").Value = 112.4972
$ws.Range("I23").Value = 0.2909830802836222
$ws.Range("J23").Value = 0.3019872047225811
$ws.Range("O23").Value = 0.1171985110386058
$ws.Range("P23").Value = 0.1175583118271966
$ws.Range("Q23").Value = 165.2132754228
$ws.Range("R23").Value = 1486.9194788052
$ws.Range("S23").Value = 0.03410278374666763
$ws.Range("T23").Value = 0.03550110598060065

$ws.Range("G24").Value = 37.49906666666667
$ws.Range("H24").Value = 112.4972
$ws.Range("I24").Value = 0.2909830802836222
$ws.Range("J24").Value = 0.3019872047225811
$ws.Range("M24").Value = 1.362560333333333
$ws.Range("N24").Value = 4.087681
$ws.Range("O24").Value = 0.03624543805965938
$ws.Range("P24").Value = 0.03635671197501131
$ws.Range("Q24").Value = 51.09474077702223
$ws.Range("R24").Value = 459.8526669932
$ws.Range("S24").Value = 0.01054680921282892
$ws.Range("T24").Value = 0.01097926182223766

$ws.Range("G25").Value = 37.49906666666667
$ws.Range("H25").Value = 112.4972
$ws.Range("I25").Value = 0.2909830802836222
$ws.Range("J25").Value = 0.3019872047225811
$ws.Range("M25").Value = 0.3451695
$ws.Range("N25").Value = 0.690339
$ws.Range("O25").Value = 0.00918184643004207
$ws.Range("P25").Value = 0.006140023203404898
$ws.Range("Q25").Value = 12.9435340918
$ws.Range("R25").Value = 77.66120455080001
$ws.Range("S25").Value = 0.002671761956904822
$ws.Range("T25").Value = 0.001854208444128033

$ws.Range("G26").Value = 37.49906666666667
$ws.Range("H26").Value = 112.4972
$ws.Range("I26").Value = 0.2909830802836222
$ws.Range("J26").Value = 0.3019872047225811
$ws.Range("M26").Value = 1.693723333333333
$ws.Range("N26").Value = 5.08117
$ws.Range("O26").Value = 0.04505469788508434
$ws.Range("P26").Value = 0.04519301633030275
$ws.Range("Q26").Value = 69.47404023088446
$ws.Range("R26").Value = 625.26636207796
$ws.Range("S26").Value = 0.01434060407033231
$ws.Range("T26").Value = 0.014928614294616
